# repull data, push all data, mean calculation
# Update column F (dSF) values for specific rows to reflect re-pulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    7  = 2
    10 = -1
    13 = 0
    16 = 2
    20 = -1
    21 = -4
    25 = 1
    26 = 4
    28 = -2
    30 = -2
    32 = 0
    35 = 0
    37 = 0
    42 = -1
    46 = -12
    48 = -2
    52 = 8
    53 = -3
    55 = -7
    57 = -2
    61 = -6
    63 = 5
    67 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
